$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hematocrit")

# --- Copy number-format styles for columns A, B, C down into the new rows ---
$ws.Range("A137:C137").Copy() | Out-Null
$ws.Range("A138:C149").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 138
$ws.Range("A138").Value = 43587
$ws.Range("B138").Value = 1
$ws.Range("C138").Value = "A"
$ws.Range("D138").Value = "20.05_8x3_1.25x7mm"
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = "M"
$ws.Range("G138").Value = 0.56000000000000005
$ws.Range("K138").Formula = "=10.0224-6.1343"
$ws.Range("L138").Value = 0.065

# Row 139
$ws.Range("A139").Value = 43587
$ws.Range("B139").Value = 2
$ws.Range("C139").Value = "A"
$ws.Range("D139").Value = "Moat2_glossy"
$ws.Range("E139").Value = 4
$ws.Range("F139").Value = "M"
$ws.Range("G139").Value = 0.56000000000000005
$ws.Range("K139").Formula = "=9.2132-6.0331"
$ws.Range("L139").Value = 0.008

# Row 140
$ws.Range("A140").Value = 43587
$ws.Range("B140").Value = 3
$ws.Range("C140").Value = "A"
$ws.Range("D140").Value = "20.05_8x3_1.25x9mm"
$ws.Range("E140").Value = 4
$ws.Range("F140").Value = "M"
$ws.Range("G140").Value = 0.56000000000000005
$ws.Range("K140").Formula = "=9.7443-6.0381"
$ws.Range("L140").Value = 0.004

# Row 141
$ws.Range("A141").Value = 43587
$ws.Range("B141").Value = 4
$ws.Range("C141").Value = "A"
$ws.Range("D141").Value = "20.05_8x3_1.25x7mm"
$ws.Range("E141").Value = 4
$ws.Range("F141").Value = "M"
$ws.Range("G141").Value = 0.56000000000000005
$ws.Range("K141").Formula = "=10.2184-6.1569"
$ws.Range("L141").Value = 0.07000000000000001

# Row 142
$ws.Range("A142").Value = 43587
$ws.Range("B142").Value = 5
$ws.Range("C142").Value = "A"
$ws.Range("D142").Value = "Moat2_glossy"
$ws.Range("E142").Value = 4
$ws.Range("F142").Value = "M"
$ws.Range("G142").Value = 0.56000000000000005
$ws.Range("K142").Formula = "=9.4637-6.0543"
$ws.Range("L142").Value = 0.005

# Row 143
$ws.Range("A143").Value = 43587
$ws.Range("B143").Value = 6
$ws.Range("C143").Value = "A"
$ws.Range("D143").Value = "20.05_8x3_1.25x9mm"
$ws.Range("E143").Value = 4
$ws.Range("F143").Value = "M"
$ws.Range("G143").Value = 0.56000000000000005
$ws.Range("K143").Formula = "=9.9504-6.0779"
$ws.Range("L143").Value = 0.045

# Row 144
$ws.Range("A144").Value = 43587
$ws.Range("B144").Value = 1
$ws.Range("C144").Value = "A"
$ws.Range("D144").Value = "20.05_8x3_1.25x7mm"
$ws.Range("E144").Value = 4
$ws.Range("F144").Value = "F"
$ws.Range("G144").Value = 0.41499999999999998
$ws.Range("I144").Interior.Pattern = -4142
$ws.Range("K144").Formula = "=10.2134-6.014"
$ws.Range("L144").Value = 0

# Row 145
$ws.Range("A145").Value = 43587
$ws.Range("B145").Value = 2
$ws.Range("C145").Value = "A"
$ws.Range("D145").Value = "20.05_8x3_1.25x9mm"
$ws.Range("E145").Value = 4
$ws.Range("F145").Value = "F"
$ws.Range("G145").Value = 0.41499999999999998
$ws.Range("I145").Interior.Pattern = -4142
$ws.Range("K145").Formula = "=10.2447-6.1587"
$ws.Range("L145").Value = 0

# Row 146
$ws.Range("A146").Value = 43587
$ws.Range("B146").Value = 3
$ws.Range("C146").Value = "A"
$ws.Range("D146").Value = "Moat2_glossy"
$ws.Range("E146").Value = 4
$ws.Range("F146").Value = "F"
$ws.Range("G146").Value = 0.41499999999999998
$ws.Range("I146").Interior.Pattern = -4142
$ws.Range("K146").Formula = "=8.9526-6.0303"
$ws.Range("L146").Value = 0

# Row 147
$ws.Range("A147").Value = 43587
$ws.Range("B147").Value = 4
$ws.Range("C147").Value = "A"
$ws.Range("D147").Value = "20.05_8x3_1.25x7mm"
$ws.Range("E147").Value = 4
$ws.Range("F147").Value = "F"
$ws.Range("G147").Value = 0.41499999999999998
$ws.Range("K147").Formula = "=10.2938-6.0305"
$ws.Range("L147").Value = 0

# Row 148
$ws.Range("A148").Value = 43587
$ws.Range("B148").Value = 5
$ws.Range("C148").Value = "A"
$ws.Range("D148").Value = "20.05_8x3_1.25x9mm"
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = "F"
$ws.Range("G148").Value = 0.41499999999999998
$ws.Range("K148").Formula = "=10.161-6.1216"
$ws.Range("L148").Value = 0

# Row 149
$ws.Range("A149").Value = 43587
$ws.Range("B149").Value = 6
$ws.Range("C149").Value = "A"
$ws.Range("D149").Value = "Moat2_glossy"
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = "F"
$ws.Range("G149").Value = 0.41499999999999998
$ws.Range("K149").Formula = "=9.5773-6.0861"
$ws.Range("L149").Value = 0


# --- Update the view: make Hematocrit the active sheet, with the newly
#     entered formula cells selected (matches the authored commit) ---
$ws.Activate()
$ws.Range("K138:K149").Select()
$excel.CutCopyMode = $false
